$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.676.10'
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").Value = '2.438.08'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.31'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.32'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.26'
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.353'
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.82'
$ws.Range("E13").Value = '  +4.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000183'
$ws.Range("E14").Value = '  +5.27%  '
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("D16").Value = '62.481.01'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").Value = '2.436.80'
$ws.Range("E17").Value = '  +1.11%  '
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.38'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.29'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("E26").Value = '  +7.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '566.55'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  +2.63%  '
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.87'
$ws.Range("E35").Value = '  +3.78%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.44'
$ws.Range("E38").Value = '  -0.97%  '
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '148.24'
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").Value = '  +2.08%  '
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '148.35'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.67'
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0537'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.48'
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.601'
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("E49").Value = '  +2.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0926'
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.08'
$ws.Range("E51").Value = '  +1.66%  '
